$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "2" = @{ "B" = 1.497172313306464; "C" = 0.1386889568703396; "D" = 0.4686512189832825; "E" = 0.1436258580032064; "G" = 0.002526051250973015; "I" = 1.326209116498454; "J" = 0.05864987886410589; "L" = 0.4556374798074643; "M" = 0.4043823545685896; "O" = 5.923380952013247 }
    "3" = @{ "B" = 1.401937764878426; "C" = 0.1225450772920453; "D" = 0.4694133419567663; "E" = 0.144675058875305; "G" = 0.002529431048235218; "I" = 1.340486853545684; "J" = 0.05769887180640509; "L" = 0.4520788122979269; "M" = 0.3885009025617379; "O" = 5.961890986998213 }
    "4" = @{ "B" = 1.343759904210515; "C" = 0.1125822551207136; "D" = 0.4700748876633298; "E" = 0.1453589704484628; "G" = 0.002531617775186718; "I" = 1.349963396191097; "J" = 0.05710957526181204; "L" = 0.4500634677193176; "M" = 0.3788683928763845; "O" = 5.988883739872108 }
    "5" = @{ "B" = 1.320128016791216; "C" = 0.1085098176845634; "D" = 0.4703932201784227; "E" = 0.1456476671964921; "G" = 0.002532537011928194; "I" = 1.354003601078613; "J" = 0.05686808978462565; "L" = 0.4492849854772913; "M" = 0.374973201884103; "O" = 6.000724816051815 }
    "6" = @{ "B" = 1.316208598266996; "C" = 0.1078328422610753; "D" = 0.4704490251889766; "E" = 0.145696209316164; "G" = 0.002532691351770542; "I" = 1.354685248625813; "J" = 0.05682791050045211; "L" = 0.4491583067772069; "M" = 0.3743282364663472; "O" = 6.002741818277627 }
    "7" = @{ "B" = 1.343440886342194; "C" = 0.1125273832159053; "D" = 0.4700789833431287; "E" = 0.1453628234156339; "G" = 0.00253163005829482; "I" = 1.350017161442267; "J" = 0.05710632392728598; "L" = 0.4500527954362283; "M" = 0.3788157385933673; "O" = 5.989040026919412 }
    "8" = @{ "B" = 1.464274757152737; "C" = 0.1331331078802123; "D" = 0.468873855719707; "E" = 0.1439793921358099; "G" = 0.002527193511197472; "I" = 1.330984677503437; "J" = 0.05832309222721932; "L" = 0.4543752995011801; "M" = 0.3988819358774123; "O" = 5.935964187834145 }
    "9" = @{ "B" = 1.703529736640974; "C" = 0.1731353089502647; "D" = 0.4680442398539668; "E" = 0.1415807871450214; "G" = 0.002519374385849036; "I" = 1.299299050380213; "J" = 0.06066623851658548; "L" = 0.4641938909685166; "M" = 0.4391645004063633; "O" = 5.858465596935815 }
    "10" = @{ "B" = 1.880660815822409; "C" = 0.202273120476292; "D" = 0.4683668841085051; "E" = 0.1400091646264334; "G" = 0.002514161227929085; "I" = 1.279461319644952; "J" = 0.06236134741941868; "L" = 0.4722214092813743; "M" = 0.469319377678957; "O" = 5.817768948563725 }
    "11" = @{ "B" = 1.961525989579968; "C" = 0.2154732280830558; "D" = 0.4687155283262854; "E" = 0.139335379716151; "G" = 0.002511903881848123; "I" = 1.271185105370833; "J" = 0.06312672693169574; "L" = 0.4760492023848713; "M" = 0.4831573130625912; "O" = 5.802789818370769 }
    "12" = @{ "B" = 1.992187652123903; "C" = 0.2204637555300906; "D" = 0.4688765255034326; "E" = 0.1390861376879311; "G" = 0.002511065409035007; "I" = 1.268158798982881; "J" = 0.06341572511260551; "L" = 0.4775239068631834; "M" = 0.4884144665082104; "O" = 5.797626406150414 }
    "13" = @{ "B" = 1.985582369871508; "C" = 0.219389317808151; "D" = 0.468840564175764; "E" = 0.1391395539542897; "G" = 0.002511245264018589; "I" = 1.2688057751839; "J" = 0.0633535214509493; "L" = 0.4772051837510389; "M" = 0.4872814915057049; "O" = 5.79871579689916 }
    "14" = @{ "B" = 1.964047754633157; "C" = 0.2158839644876309; "D" = 0.4687281933131118; "E" = 0.1393147561351427; "G" = 0.002511834573236602; "I" = 1.270933969695186; "J" = 0.0631505197652622; "L" = 0.4761700227526973; "M" = 0.4835894825517357; "O" = 5.802354819157074 }
    "15" = @{ "B" = 1.950862315437632; "C" = 0.2137357784459084; "D" = 0.4686631344297325; "E" = 0.1394228412647243; "G" = 0.002512197667260188; "I" = 1.27225158244071; "J" = 0.06302606641277819; "L" = 0.4755392355475863; "M" = 0.4813302298238682; "O" = 5.804650113991244 }
    "16" = @{ "B" = 1.875381736450322; "C" = 0.2014093434533493; "D" = 0.4683481567934393; "E" = 0.1400540251340843; "G" = 0.002514311041077251; "I" = 1.280017255149005; "J" = 0.06231121183704147; "L" = 0.4719747867363111; "M" = 0.468417433785433; "O" = 5.818819040520083 }
    "17" = @{ "B" = 1.829149340735569; "C" = 0.1938333012908515; "D" = 0.4682065926550933; "E" = 0.1404517678191981; "G" = 0.002515636706568676; "I" = 1.284972945200892; "J" = 0.06187119601627344; "L" = 0.469833120104127; "M" = 0.4605264741790904; "O" = 5.828416836464214 }
    "18" = @{ "B" = 1.802584828362114; "C" = 0.1894706122266712; "D" = 0.4681441731429743; "E" = 0.1406844132306162; "G" = 0.002516409942882021; "I" = 1.287893730125752; "J" = 0.06161757142303514; "L" = 0.4686178659237044; "M" = 0.455999145315289; "O" = 5.834269829063373 }
    "19" = @{ "B" = 1.793595249376438; "C" = 0.1879926017485047; "D" = 0.4681263047018405; "E" = 0.1407638487644292; "G" = 0.002516673596184165; "I" = 1.288894745109609; "J" = 0.0615316062010649; "L" = 0.4682092520577896; "M" = 0.4544682260630424; "O" = 5.836308654574566 }
    "20" = @{ "B" = 1.834068060659263; "C" = 0.1946403182460585; "D" = 0.4682196959166447; "E" = 0.1404090265406213; "G" = 0.002515494475349433; "I" = 1.284438116027729; "J" = 0.06191809232191758; "L" = 0.4700593894469876; "M" = 0.4613653084105209; "O" = 5.827360707526339 }
    "21" = @{ "B" = 1.970371925859752; "C" = 0.2169137921051458; "D" = 0.4687604134784351; "E" = 0.139263134825713; "G" = 0.002511661036118938; "I" = 1.270305942830788; "J" = 0.06321016900527354; "L" = 0.4764733915293817; "M" = 0.4846734551367931; "O" = 5.801272134408038 }
    "22" = @{ "B" = 2.059685521263873; "C" = 0.2314236920702513; "D" = 0.4692826558912202; "E" = 0.1385486456531408; "G" = 0.002509250840626365; "I" = 1.261697671594604; "J" = 0.06404974761555238; "L" = 0.4808121599828752; "M" = 0.5000058003674468; "O" = 5.787188080546514 }
    "23" = @{ "B" = 2.01199653184625; "C" = 0.2236838550255129; "D" = 0.4689884936263269; "E" = 0.1389268367957577; "G" = 0.002510528524373282; "I" = 1.266234565138966; "J" = 0.06360209738528511; "L" = 0.4784830775985114; "M" = 0.4918136597623004; "O" = 5.794433360977109 }
    "24" = @{ "B" = 1.831844260206253; "C" = 0.1942754880534494; "D" = 0.4682137128466479; "E" = 0.140428337481195; "G" = 0.002515558743474687; "I" = 1.284679688983211; "J" = 0.06189689253718456; "L" = 0.4699570431731814; "M" = 0.4609860424998473; "O" = 5.827837139654179 }
    "25" = @{ "B" = 1.638564093208515; "C" = 0.1623575521518319; "D" = 0.4681047755509127; "E" = 0.1421961263259695; "G" = 0.002521395925878219; "I" = 1.307266832883876; "J" = 0.060036975618317; "L" = 0.461394473909877; "M" = 0.4641938909685166; "O" = 5.876581427243224 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $data[$row][$col]
    }
}
